# "Hortaliza, Terminal La Palmera de La Serena - Jengibre"
# Weekly refresh: a new weekly observation is inserted as row 157,
# pushing the existing rows 157-200 down to 158-201 (dimension grows
# from A1:R200 to A1:R201).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 157, shifting rows 157:200 down to 158:201.
$ws.Rows.Item(157).Insert()

# Populate the newly inserted row 157 with the new weekly record.
$ws.Cells.Item(157, 1).Value = 8
$ws.Cells.Item(157, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(157, 3).Value = "Coquimbo"
$ws.Cells.Item(157, 4).Value = 45244
$ws.Cells.Item(157, 5).Value = 4
$ws.Cells.Item(157, 6).Value = 100114007
$ws.Cells.Item(157, 7).Value = "Jengibre"
$ws.Cells.Item(157, 8).Value = "Sin especificar"
$ws.Cells.Item(157, 9).Value = "Primera"
$ws.Cells.Item(157, 10).Value = 400
$ws.Cells.Item(157, 11).Value = 24000
$ws.Cells.Item(157, 12).Value = 25000
$ws.Cells.Item(157, 13).Value = 24500
$ws.Cells.Item(157, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(157, 15).Value = "Perú"
$ws.Cells.Item(157, 16).Value = 1885
$ws.Cells.Item(157, 17).Value = 13
$ws.Cells.Item(157, 18).Value = "Hortaliza"
